$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "BBBKN3DD&E"
$ws.Range("B2").Value = "SYMBOL_200"

$ws.Range("E7").Select()
